$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("Recorded By") holds entries like "dnasr281@gmail.com, System".
# Swap the order of the two recorder names to "System, dnasr281@gmail.com"
# wherever that exact combination appears.
$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$colG = $ws.Range("G1:G$lastRow")
[void]$colG.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com", 1, 1, $false, $false, $false, $false)
